{"js": "// Map of old text -> new text for every <w:t> run that changes.\nconst replacements = [\n  [\"2025-01-02 Thursday\", \"2025-01-03 Friday\"],\n  [\"35\u00d746=1610\", \"89\u00d712=1068\"],\n  [\"92\u00d769=6348\", \"43\u00d731=1333\"],\n  [\"28\u00d715=420\", \"13\u00d785=1105\"],\n  [\"19\u00d769=1311\", \"57\u00d742=2394\"],\n  [\"51\u00d776=3876\", \"73\u00d737=2701\"],\n  [\"47\u00d788=4136\", \"71\u00d770=4970\"],\n  [\"39\u00d740=1560\", \"52\u00d749=2548\"],\n  [\"17\u00d787=1479\", \"74\u00d776=5624\"],\n  [\"39\u00d746=1794\", \"76\u00d787=6612\"],\n  [\"58\u00d786=4988\", \"42\u00d717=714\"],\n  [\"79\u00d715=1185\", \"13\u00d734=442\"],\n  [\"94\u00d797=9118\", \"60\u00d779=4740\"],\n  [\"42\u00d721=882\", \"29\u00d739=1131\"],\n  [\"72\u00d728=2016\", \"18\u00d772=1296\"],\n  [\"51\u00d769=3519\", \"46\u00d712=552\"],\n  [\"55\u00d718=990\", \"22\u00d743=946\"],\n  [\"52\u00d744=2288\", \"61\u00d798=5978\"],\n  [\"35\u00d775=2625\", \"12\u00d771=852\"],\n  [\"41\u00d791=3731\", \"29\u00d757=1653\"],\n  [\"57\u00d773=4161\", \"99\u00d712=1188\"],\n  [\"48\u00d767=3216\", \"35\u00d725=875\"],\n  [\"38\u00d760=2280\", \"69\u00d742=2898\"],\n  [\"95\u00d783=7885\", \"82\u00d798=8036\"],\n  [\"80\u00d724=1920\", \"11\u00d793=1023\"],\n  [\"77\u00d772=5544\", \"48\u00d744=2112\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old=\"2025-01-02 Thursday\"; new=\"2025-01-03 Friday\"},\n  @{old=\"35\u00d746=1610\"; new=\"89\u00d712=1068\"},\n  @{old=\"92\u00d769=6348\"; new=\"43\u00d731=1333\"},\n  @{old=\"28\u00d715=420\"; new=\"13\u00d785=1105\"},\n  @{old=\"19\u00d769=1311\"; new=\"57\u00d742=2394\"},\n  @{old=\"51\u00d776=3876\"; new=\"73\u00d737=2701\"},\n  @{old=\"47\u00d788=4136\"; new=\"71\u00d770=4970\"},\n  @{old=\"39\u00d740=1560\"; new=\"52\u00d749=2548\"},\n  @{old=\"17\u00d787=1479\"; new=\"74\u00d776=5624\"},\n  @{old=\"39\u00d746=1794\"; new=\"76\u00d787=6612\"},\n  @{old=\"58\u00d786=4988\"; new=\"42\u00d717=714\"},\n  @{old=\"79\u00d715=1185\"; new=\"13\u00d734=442\"},\n  @{old=\"94\u00d797=9118\"; new=\"60\u00d779=4740\"},\n  @{old=\"42\u00d721=882\"; new=\"29\u00d739=1131\"},\n  @{old=\"72\u00d728=2016\"; new=\"18\u00d772=1296\"},\n  @{old=\"51\u00d769=3519\"; new=\"46\u00d712=552\"},\n  @{old=\"55\u00d718=990\"; new=\"22\u00d743=946\"},\n  @{old=\"52\u00d744=2288\"; new=\"61\u00d798=5978\"},\n  @{old=\"35\u00d775=2625\"; new=\"12\u00d771=852\"},\n  @{old=\"41\u00d791=3731\"; new=\"29\u00d757=1653\"},\n  @{old=\"57\u00d773=4161\"; new=\"99\u00d712=1188\"},\n  @{old=\"48\u00d767=3216\"; new=\"35\u00d725=875\"},\n  @{old=\"38\u00d760=2280\"; new=\"69\u00d742=2898\"},\n  @{old=\"95\u00d783=7885\"; new=\"82\u00d798=8036\"},\n  @{old=\"80\u00d724=1920\"; new=\"11\u00d793=1023\"},\n  @{old=\"77\u00d772=5544\"; new=\"48\u00d744=2112\"}\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $p.old\n  $find.Replacement.Text = $p.new\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
